$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 118
$ws.Range("H3").Value = 159
$ws.Range("C6").Value = 493
$ws.Range("E6").Value = 492
$ws.Range("F6").Value = 560
$ws.Range("G6").Value = 441
$ws.Range("H6").Value = 460
$ws.Range("I6").Value = 511
$ws.Range("C7").Value = 651
$ws.Range("E7").Value = 727
$ws.Range("F7").Value = 810
$ws.Range("G7").Value = 677
$ws.Range("H7").Value = 746
$ws.Range("I7").Value = 852

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("H3").Value = 8
$ws.Range("H7").Value = 38

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("G6").Value = 31
$ws.Range("G7").Value = 47

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F8").Value = 53
$ws.Range("C11").Value = 7
$ws.Range("I19").Value = 24
$ws.Range("I21").Value = 17
$ws.Range("G28").Value = 47
$ws.Range("H36").Value = 38
$ws.Range("H53").Value = 113
$ws.Range("E54").Value = 6
$ws.Range("I61").Value = 4
$ws.Range("E65").Value = 17
$ws.Range("C98").Value = 651
$ws.Range("E98").Value = 727
$ws.Range("F98").Value = 810
$ws.Range("G98").Value = 677
$ws.Range("H98").Value = 746
$ws.Range("I98").Value = 852

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H6").Value = 75
$ws.Range("H7").Value = 113

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("E5").Value = 14
$ws.Range("E6").Value = 17

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 7

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 2
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 6

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F5").Value = 38
$ws.Range("F6").Value = 53

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 17
